$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOV-2020")

# --- Fix row 12 (was a duplicated "10" row referencing the Hayaai-site text) ---
$ws.Range("A12").Value = 11
$ws.Range("C12").Value = "B2C & Sonia"
$ws.Range("D12").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing, Retesting on B2C/B2B app and Sonia application"

# --- Rows 13 & 14: same pattern as (the now-fixed) row 12 ---
$ws.Range("A12:G12").Copy($ws.Range("A13:G13"))
$ws.Rows.Item(13).RowHeight = 30
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 44147

$ws.Range("A12:G12").Copy($ws.Range("A14:G14"))
$ws.Rows.Item(14).RowHeight = 30
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 44148

# --- Rows 15 & 16: "Week off" rows, same pattern as row 8/9 ---
$ws.Range("A8:G8").Copy($ws.Range("A15:G15"))
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 44149

$ws.Range("A8:G8").Copy($ws.Range("A16:G16"))
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 44150

# --- Row 17: new data row (B2C & B2B). Copy cell-by-cell so the destination
# row does NOT inherit the legacy row-level customFormat flag. ---
$ws.Range("A12").Copy($ws.Range("A17"))
$ws.Range("B12").Copy($ws.Range("B17"))
$ws.Range("G12").Copy($ws.Range("C17"))
$ws.Range("D12").Copy($ws.Range("D17"))
$ws.Range("E12").Copy($ws.Range("E17"))
$ws.Range("F12").Copy($ws.Range("F17"))
$ws.Range("G12").Copy($ws.Range("G17"))
$ws.Rows.Item(17).RowHeight = 30
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 44151
$ws.Range("C17").Value = "B2C & B2B"
$ws.Range("D17").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing, Retesting on B2C/B2B app"

# --- Row 18: mostly empty trailer row ---
$ws.Range("A12").Copy($ws.Range("A18"))
$ws.Range("B12").Copy($ws.Range("B18"))
$ws.Range("G12").Copy($ws.Range("C18"))
$ws.Range("D12").Copy($ws.Range("D18"))
$ws.Range("G12").Copy($ws.Range("E18"))
$ws.Range("G12").Copy($ws.Range("F18"))
$ws.Range("G12").Copy($ws.Range("G18"))
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 44152
$ws.Range("D18").ClearContents()

# --- Move the legend block from B20:C23/C19 into I2:I5/J1:J5 ---
$ws.Range("C19").Copy($ws.Range("J1"))
$ws.Range("B20").Copy($ws.Range("I2"))
$ws.Range("C20").Copy($ws.Range("J2"))
$ws.Range("B21").Copy($ws.Range("I3"))
$ws.Range("C21").Copy($ws.Range("J3"))
$ws.Range("B22").Copy($ws.Range("I4"))
$ws.Range("C22").Copy($ws.Range("J4"))
$ws.Range("B23").Copy($ws.Range("I5"))
$ws.Range("C23").Copy($ws.Range("J5"))
$ws.Range("B19:C23").Clear()

# --- View state: scroll down a bit, select I22 ---
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("I22").Select()
